$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$TAB = [char]9
$NL  = [char]10

# --- B2 (CasesTab query): drop the leading tab before the last line and
#     append the new "order By ... LIMIT 100" clause. ---
$oldTail2 = $TAB + "demo.survival_time AS ``Survival (days)``"
$newTail2 = "demo.survival_time AS ``Survival (days)``"
$b2 = $ws.Range("B2").Value2
$b2 = $b2.Replace($oldTail2, $newTail2)
$b2 = $b2 + $NL + " order By ss.study_subject_id ASC LIMIT 100"
$ws.Range("B2").Value = $b2

# --- B3 (SamplesTab query): just append the new "order By ... LIMIT 100"
#     clause on a new line. ---
$b3 = $ws.Range("B3").Value2
$b3 = $b3 + $NL + " order By samp.sample_id ASC LIMIT 100"
$ws.Range("B3").Value = $b3

# --- B4 (FilesTab query): replace the trailing "    order by f.file_name"
#     with "     order By f.file_name ASC LIMIT 100". ---
$oldTail4 = "    order by f.file_name"
$newTail4 = "     order By f.file_name ASC LIMIT 100"
$b4 = $ws.Range("B4").Value2
$b4 = $b4.Replace($oldTail4, $newTail4)
$ws.Range("B4").Value = $b4

# --- Row heights grow by one wrapped line (14.4pt) for rows 2 and 3; row 4
#     keeps its original height since its text only grew on the existing
#     last line (no new wrapped line introduced). ---
$ws.Rows.Item(2).RowHeight = 403.2
$ws.Rows.Item(3).RowHeight = 360
